# "Generate Report for Handoff"
#
# The localization status report moved from "In Translation" to
# "Ready for handoff", and the two "Latest ... Datetime" timestamps that
# are stamped at report-generation time were refreshed a few seconds later
# (19:18:59 -> 19:19:38, and 19:18:54 -> 19:19:33). The two status columns
# on the per-locale tabs also grew wider to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---------------
$overview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$overview.Range("F2").Value = "Ready for handoff"   # de-de status
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed timestamps -------------------------------------------------
# Overview!G2 and de-de!H2 share the "Latest HO Xliff Generate Date" /
# "Latest Handoff Datetime" value for the de-de handoff.
$overview.Range("G2").Value = "2016-09-06 19:19:38"
$dede.Range("H2").Value     = "2016-09-06 19:19:38"

# zh-cn!H2 "Latest Handoff Datetime" for the zh-cn handoff.
$zhcn.Range("H2").Value = "2016-09-06 19:19:33"

# --- Widen the Status columns to fit "Ready for handoff" ------------------
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332  # E: zh-cn status
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332  # F: de-de status
$zhcn.Columns.Item(3).ColumnWidth     = 16.333333333333332  # C: Status
$dede.Columns.Item(3).ColumnWidth     = 16.333333333333332  # C: Status
